$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old rows 13-24 down to 14-25)
$ws.Rows.Item(13).Insert()

# Row 10 (Objetivos): replace short text with the long Portuguese objectives text
$ws.Range("B10").Value = 'Gerais - Apresentar e Ensinar conceitos de Química Orgânica como instrumentos importantes para a compreensão de estratégias e operações industriais e tecnológicas. Abordar problemáticas sociais e ambientais com as quais a engenharia química está relacionada, tornando-os dessa forma, aptos a exercerem a função de Engenheiro Químico, e realizarem as mudanças que se façam necessárias.Específicos – Compreender e descrever o mecanismo das reações orgânicas e a sua importância para o aprimoramento e desenvolvimento de processos industriais sintéticos e de etapas de formulação. Aprofundar o conceito de estrutura-reatividade e propriedades dos materiais.'
$ws.Range("C10").Value = 'Gerais - Apresentar e Ensinar conceitos de Química Orgânica como instrumentos importantes para a compreensão de estratégias e operações industriais e tecnológicas. Abordar problemáticas sociais e ambientais com as quais a engenharia química está relacionada, tornando-os dessa forma, aptos a exercerem a função de Engenheiro Químico, e realizarem as mudanças que se façam necessárias.Específicos – Compreender e descrever o mecanismo das reações orgânicas e a sua importância para o aprimoramento e desenvolvimento de processos industriais sintéticos e de etapas de formulação. Aprofundar o conceito de estrutura-reatividade e propriedades dos materiais.'

# New row 13 (Docentes responsáveis value): teacher name
$ws.Range("B13").Value = '210064 - Eduardo Rezende Triboni'
$ws.Range("C13").Value = '210064 - Eduardo Rezende Triboni'

# Row 14 (Programa resumido): long Portuguese summary replacing "Semestral"
$ws.Range("B14").Value = 'Propriedade gerais dos compostos orgânicos. Estrutura, métodos de obtenção, propriedades físicas, reações dos hidrocarbonetos alifáticos e aromáticos, haletos orgânicos, álcoois e características estruturaiscomo estereoquímica e a relação estrutura-reatividade.'
$ws.Range("C14").Value = 'Propriedade gerais dos compostos orgânicos. Estrutura, métodos de obtenção, propriedades físicas, reações dos hidrocarbonetos alifáticos e aromáticos, haletos orgânicos, álcoois e características estruturaiscomo estereoquímica e a relação estrutura-reatividade.'

# Row 16 (Programa): long Portuguese syllabus replacing date placeholder
$ws.Range("B16").Value = '1.Teoria de Bronsted e de Lewis e acidez de compostos orgânicos2.Alcanos - Processos de obtenção, Propriedades físicas, Análise Conformacional. Reação de Substituição Radicalar. 3.Isomeria Constitucional e Isomeria Espacial (Estereoquímica). Quiralidade, Nomenclatura R/S, classificação de estereoisômeros. Polarímetro e Técnicas de  Resolução de Isômeros Espaciais.4.Haletos de Alquila – Substituição Nucleofílica, SN1, SN2, E1, E2. 5.Alcenos, Alcadienos e Alcinos – Propriedades físicas e químicas. Reação de adição eletrofílica (hidroalogenação, Hidratação, Halogenação, Diels-Alder, Redução-Oxidação). Adição conjugada em dienos (produto termodinâmico e cinético) 6. Fundamentos de RMN, Infra-vermelho, Ultra-violeta e Fluorescencia 7.Compostos aromáticos – Propriedades físicas dos aromáticos. Reações de Substituição Eletrofílica Aromática. Efeito de Grupos Substituintes. Reação de Substituição Nucleofílica.8.Álcoois e Éteres – Propriedades físicas, reações e mecanismos.'
$ws.Range("C16").Value = '1.Teoria de Bronsted e de Lewis e acidez de compostos orgânicos2.Alcanos - Processos de obtenção, Propriedades físicas, Análise Conformacional. Reação de Substituição Radicalar. 3.Isomeria Constitucional e Isomeria Espacial (Estereoquímica). Quiralidade, Nomenclatura R/S, classificação de estereoisômeros. Polarímetro e Técnicas de  Resolução de Isômeros Espaciais.4.Haletos de Alquila – Substituição Nucleofílica, SN1, SN2, E1, E2. 5.Alcenos, Alcadienos e Alcinos – Propriedades físicas e químicas. Reação de adição eletrofílica (hidroalogenação, Hidratação, Halogenação, Diels-Alder, Redução-Oxidação). Adição conjugada em dienos (produto termodinâmico e cinético) 6. Fundamentos de RMN, Infra-vermelho, Ultra-violeta e Fluorescencia 7.Compostos aromáticos – Propriedades físicas dos aromáticos. Reações de Substituição Eletrofílica Aromática. Efeito de Grupos Substituintes. Reação de Substituição Nucleofílica.8.Álcoois e Éteres – Propriedades físicas, reações e mecanismos.'

# Row 19 (Método): evaluation method text
$ws.Range("B19").Value = 'Duas provas teóricas e ao longo do semestre letivoAos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno.'
$ws.Range("C19").Value = 'Duas provas teóricas e ao longo do semestre letivoAos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno.'

# Row 20 (Critério): grading criteria formula text
$ws.Range("B20").Value = 'A média final (M) será calculada pela expressão: M = (P1 + P2)/2'
$ws.Range("C20").Value = 'A média final (M) será calculada pela expressão: M = (P1 + P2)/2'

# Row 21 (Norma de recuperação): recovery norm text
$ws.Range("B21").Value = 'Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno.'
$ws.Range("C21").Value = 'Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno.'

# Row 22 (Bibliografia): full bibliography list
$ws.Range("B22").Value = 'BRESLOW, R. Questões e Exercícios de Química Orgânica. São Paulo: Makrons Books Editora, 1996. 
BRUICE, P. Y. Química Orgânica, vol 1 e 2, São Paulo: Editora Pearson Prentice Hall, 2006. 
HENDRIKSON, James B.; CRAM, Donald J. Mecanismos de Reações Orgânicas. São Paulo: Livraria Editora, 1966.
MCMURRY, John. Química Orgânica. São Paulo: Editora Pioneira Thomson Leraning, 2005.
SOLOMONS, T.W.G; FRYHLE, Graig. Química Orgânica. Rio de Janeiro: Livros Técnicos e Científicos Editora, 2001.
MORRISON, R.; BOYD, R. Química Orgânica. São Paulo: Editora Calouste Gulbenkian, 2008.'
$ws.Range("C22").Value = 'BRESLOW, R. Questões e Exercícios de Química Orgânica. São Paulo: Makrons Books Editora, 1996. 
BRUICE, P. Y. Química Orgânica, vol 1 e 2, São Paulo: Editora Pearson Prentice Hall, 2006. 
HENDRIKSON, James B.; CRAM, Donald J. Mecanismos de Reações Orgânicas. São Paulo: Livraria Editora, 1966.
MCMURRY, John. Química Orgânica. São Paulo: Editora Pioneira Thomson Leraning, 2005.
SOLOMONS, T.W.G; FRYHLE, Graig. Química Orgânica. Rio de Janeiro: Livros Técnicos e Científicos Editora, 2001.
MORRISON, R.; BOYD, R. Química Orgânica. São Paulo: Editora Calouste Gulbenkian, 2008.'

# Column layout fix: column A should only use the 30.71 width definition
# (previously it spanned columns A:B before B got its own 60.71 width override)
$ws.Columns.Item(1).ColumnWidth = 29.83
